# Auto-generated Excel COM-interop script applying the scheduled-runner price update.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for the rows
# affected by the refreshed market data, across all 8 Job sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 117.6
$ws.Range("I12").Value = 117.6
$ws.Range("K12").Value = 117.6
$ws.Range("M12").Value = 52.40000000000001

$ws.Range("H33").Value = 475
$ws.Range("I33").Value = 350
$ws.Range("J33").Value = 600
$ws.Range("K33").Value = 350
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = -121
$ws.Range("N33").Value = -1058

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()  # was -5374

$ws.Range("H112").Value = 3469.5833
$ws.Range("J112").Value = 3999.1
$ws.Range("L112").Value = 11997.3
$ws.Range("N112").Value = -14213.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1034.8334
$ws.Range("I45").Value = 1034.8334
$ws.Range("K45").Value = 1034.8334
$ws.Range("M45").Value = -657.8334

$ws.Range("H74").Value = 2543.5833
$ws.Range("I74").Value = 1911.625
$ws.Range("J74").Value = 3807.5
$ws.Range("K74").Value = 1911.625
$ws.Range("L74").Value = 3807.5
$ws.Range("M74").Value = -1037.625
$ws.Range("N74").Value = -5555.5

$ws.Range("H77").Value = 2543.5833
$ws.Range("I77").Value = 1911.625
$ws.Range("J77").Value = 3807.5
$ws.Range("K77").Value = 9558.125
$ws.Range("L77").Value = 19037.5
$ws.Range("M77").Value = -5190.125
$ws.Range("N77").Value = -27773.5

$ws.Range("H110").Value = 1766
$ws.Range("I110").Value = 1957.5
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 1957.5
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 87.5
$ws.Range("N110").Value = -5090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7833.222
$ws.Range("I94").Value = 8675
$ws.Range("J94").Value = 1099
$ws.Range("K94").Value = 8675
$ws.Range("L94").Value = 1099
$ws.Range("M94").Value = -8224
$ws.Range("N94").Value = -2001

$ws.Range("H99").Value = 2336.2
$ws.Range("I99").Value = 1593.75
$ws.Range("J99").Value = 3184.7144
$ws.Range("K99").Value = 1593.75
$ws.Range("L99").Value = 3184.7144
$ws.Range("M99").Value = -95.75
$ws.Range("N99").Value = -6180.7144

$ws.Range("H105").Value = 2670.2307
$ws.Range("I105").Value = 2710.3635
$ws.Range("J105").Value = 2449.5
$ws.Range("K105").Value = 2710.3635
$ws.Range("L105").Value = 2449.5
$ws.Range("M105").Value = -963.3634999999999
$ws.Range("N105").Value = -5943.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1636.2106
$ws.Range("I16").Value = 1141.2858
$ws.Range("K16").Value = 1141.2858
$ws.Range("M16").Value = -854.2858000000001

$ws.Range("H31").Value = 3085.1538
$ws.Range("I31").Value = 2419
$ws.Range("K31").Value = 2419
$ws.Range("M31").Value = -2124

$ws.Range("H34").Value = 3085.1538
$ws.Range("I34").Value = 2419
$ws.Range("K34").Value = 2419
$ws.Range("M34").Value = -2217

$ws.Range("H86").Value = 19966.334
$ws.Range("I86").Value = 19966.334
$ws.Range("K86").Value = 19966.334
$ws.Range("M86").Value = -18843.334

$ws.Range("H88").Value = 11921.333
$ws.Range("J88").Value = 11921.333
$ws.Range("L88").Value = 11921.333
$ws.Range("N88").Value = -12733.333

$ws.Range("H89").Value = 19966.334
$ws.Range("I89").Value = 19966.334
$ws.Range("K89").Value = 99831.67
$ws.Range("M89").Value = -94215.67

$ws.Range("H91").Value = 11921.333
$ws.Range("J91").Value = 11921.333
$ws.Range("L91").Value = 11921.333
$ws.Range("N91").Value = -14729.333

$ws.Range("H99").Value = 3848.75
$ws.Range("I99").Value = 1995
$ws.Range("J99").Value = 4113.5713
$ws.Range("K99").Value = 1995
$ws.Range("L99").Value = 4113.5713
$ws.Range("M99").Value = -497
$ws.Range("N99").Value = -7109.5713

$ws.Range("H107").Value = 905
$ws.Range("I107").Value = 522.8333
$ws.Range("K107").Value = 522.8333
$ws.Range("M107").Value = 1397.1667

$ws.Range("H113").Value = 1636.2106
$ws.Range("I113").Value = 1141.2858
$ws.Range("K113").Value = 1141.2858
$ws.Range("M113").Value = 1028.7142

$ws.Range("H126").Value = 3848.75
$ws.Range("I126").Value = 1995
$ws.Range("J126").Value = 4113.5713
$ws.Range("K126").Value = 5985
$ws.Range("L126").Value = 12340.7139
$ws.Range("M126").Value = -3515
$ws.Range("N126").Value = -17280.7139

$ws.Range("H132").Value = 2752
$ws.Range("I132").Value = 1534.8
$ws.Range("J132").Value = 3766.3333
$ws.Range("K132").Value = 4604.4
$ws.Range("L132").Value = 11298.9999
$ws.Range("M132").Value = -2074.4
$ws.Range("N132").Value = -16358.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 23904842
$ws.Range("I4").Value = 27500102
$ws.Range("K4").Value = 82500306
$ws.Range("M4").Value = -82500194

$ws.Range("H74").Value = 21000
$ws.Range("J74").Value = 21000
$ws.Range("L74").Value = 63000
$ws.Range("N74").Value = -65122

$ws.Range("H77").Value = 21000
$ws.Range("J77").Value = 21000
$ws.Range("L77").Value = 189000
$ws.Range("N77").Value = -199608

$ws.Range("H107").Value = 850.3333
$ws.Range("J107").Value = 1333.3334
$ws.Range("L107").Value = 4000.0002
$ws.Range("N107").Value = -7840.0002

$ws.Range("H130").Value = 2583.3333
$ws.Range("J130").Value = 3163.3333
$ws.Range("L130").Value = 9489.999899999999
$ws.Range("N130").Value = -19529.9999

$ws.Range("H131").Value = 2603
$ws.Range("I131").Value = 1549.8
$ws.Range("J131").Value = 3081.7273
$ws.Range("K131").Value = 4649.4
$ws.Range("L131").Value = 9245.1819
$ws.Range("M131").Value = 390.6000000000004
$ws.Range("N131").Value = -19325.1819

$ws.Range("H134").Value = 14209.25
$ws.Range("I134").Value = 1700
$ws.Range("J134").Value = 18379
$ws.Range("K134").Value = 5100
$ws.Range("L134").Value = 55137
$ws.Range("M134").Value = -30
$ws.Range("N134").Value = -65277

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 107.35294
$ws.Range("I2").Value = 135.6923
$ws.Range("J2").Value = 15.25
$ws.Range("K2").Value = 135.6923
$ws.Range("L2").Value = 15.25
$ws.Range("M2").Value = -22.69229999999999
$ws.Range("N2").Value = -241.25

$ws.Range("H52").Value = 49500
$ws.Range("I52").Value = 49500
$ws.Range("K52").Value = 49500
$ws.Range("M52").Value = -49241

$ws.Range("H113").Value = 1176.2222
$ws.Range("I113").Value = 1132
$ws.Range("K113").Value = 1132
$ws.Range("M113").Value = 1038

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3168.1667
$ws.Range("J68").Value = 3168.1667
$ws.Range("L68").Value = 3168.1667
$ws.Range("N68").Value = -4666.1667

$ws.Range("H71").Value = 3168.1667
$ws.Range("J71").Value = 3168.1667
$ws.Range("L71").Value = 15840.8335
$ws.Range("N71").Value = -23328.8335

$ws.Range("H100").Value = 4612.375
$ws.Range("I100").Value = 3633
$ws.Range("K100").Value = 3633
$ws.Range("M100").Value = -3092

$ws.Range("H129").Value = 66494
$ws.Range("J129").Value = 66494
$ws.Range("L129").Value = 66494
$ws.Range("N129").Value = -76494

$ws.Range("H139").Value = 89420
$ws.Range("I139").Value = 89420
$ws.Range("K139").Value = 89420
$ws.Range("M139").Value = -84280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()  # was -46039

$ws.Range("H129").Value = 71999.75
$ws.Range("J129").Value = 71999.75
$ws.Range("L129").Value = 71999.75
$ws.Range("N129").Value = -81999.75

Write-Output "Applied scheduled market-price refresh across all job sheets."